$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.117.99"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "1.794.13"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5345"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3762"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.095"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("E13").Value = "  -2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.104"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.784.66"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.219"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001058"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06456"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9988"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.902"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").Value = "28.147.73"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.095"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").Value = "1.988.94"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.284"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.117"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1047"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.655"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.572"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.45%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06524"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("B36").Value = "Algorand"
$ws.Range("C36").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2256"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02283"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.036"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.453"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.450"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.45%  "
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.176"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9989"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.674"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5777"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  +3.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.928"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06816"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.20%  "
